$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added to the dataset. It is inserted as row 177,
# pushing the previously-existing rows 177..224 down to 178..225.
$ws.Rows.Item(177).Insert()

$ws.Range("A177").Value = 10
$ws.Range("B177").Value = "Vega Modelo de Temuco"
$ws.Range("C177").Value = "La Araucanía"
$ws.Range("D177").Value = 45204
$ws.Range("E177").Value = 9
$ws.Range("F177").Value = "Fruta"
$ws.Range("G177").Value = 100107
$ws.Range("H177").Value = "Otros"
$ws.Range("I177").Value = 100107002
$ws.Range("J177").Value = "Chirimoya"
$ws.Range("K177").Value = "Cultivar IV Región"
$ws.Range("L177").Value = "Primera"
$ws.Range("M177").Value = 400
$ws.Range("N177").Value = 2500
$ws.Range("O177").Value = 2600
$ws.Range("P177").Value = 2562
$ws.Range("Q177").Value = "`$/kilo (en caja de 15 kilos)"
$ws.Range("R177").Value = "Provincia del Elquí"
$ws.Range("S177").Value = 2562
$ws.Range("T177").Value = 1

# Match the date cell style used by the other rows in column D.
$ws.Range("D177").NumberFormat = $ws.Range("D178").NumberFormat
